$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 220.2
$ws.Range("I8").Value = 25.5
$ws.Range("J8").Value = 999
$ws.Range("K8").Value = 76.5
$ws.Range("L8").Value = 2997
$ws.Range("M8").Value = 62.5
$ws.Range("N8").Value = -3275
$ws.Range("H11").Value = 152
$ws.Range("I11").Value = 152
$ws.Range("K11").Value = 152
$ws.Range("M11").Value = -12
$ws.Range("H19").Value = 1165.2142
$ws.Range("I19").Value = 1124.6666
$ws.Range("K19").Value = 1124.6666
$ws.Range("M19").Value = -949.6666
$ws.Range("H32").Value = 40187.25
$ws.Range("I32").Value = 68624.75
$ws.Range("K32").Value = 68624.75
$ws.Range("M32").Value = -68298.75
$ws.Range("H33").Value = 3380015.5
$ws.Range("I33").Value = 5405804.5
$ws.Range("J33").Value = 3700
$ws.Range("K33").Value = 5405804.5
$ws.Range("L33").Value = 3700
$ws.Range("M33").Value = -5405575.5
$ws.Range("N33").Value = -4158
$ws.Range("H38").Value = 1284.9546
$ws.Range("I38").Value = 373
$ws.Range("J38").Value = 3239.1428
$ws.Range("K38").Value = 1119
$ws.Range("L38").Value = 9717.428400000001
$ws.Range("M38").Value = -747
$ws.Range("N38").Value = -10461.4284
$ws.Range("H39").Value = 760.5
$ws.Range("I39").Value = 497.7143
$ws.Range("K39").Value = 1493.1429
$ws.Range("M39").Value = -1197.1429
$ws.Range("H43").Value = 3933.2
$ws.Range("I43").Value = 3891.5
$ws.Range("J43").Value = 4100
$ws.Range("K43").Value = 3891.5
$ws.Range("L43").Value = 4100
$ws.Range("M43").Value = -3822.5
$ws.Range("N43").Value = -4238
$ws.Range("H53").Value = 248.5
$ws.Range("I53").Value = 267.63635
$ws.Range("K53").Value = 267.63635
$ws.Range("M53").Value = 369.36365
$ws.Range("H64").Value = 3600
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3352
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3600
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2742
$ws.Range("N67").ClearContents()
$ws.Range("H86").Value = 4930.6665
$ws.Range("I86").Value = 2008.6666
$ws.Range("J86").Value = 9313.666999999999
$ws.Range("K86").Value = 2008.6666
$ws.Range("L86").Value = 9313.666999999999
$ws.Range("M86").Value = -885.6666
$ws.Range("N86").Value = -11559.667
$ws.Range("H88").Value = 5494.0835
$ws.Range("J88").Value = 5286.5713
$ws.Range("L88").Value = 5286.5713
$ws.Range("N88").Value = -6098.5713
$ws.Range("H89").Value = 4930.6665
$ws.Range("I89").Value = 2008.6666
$ws.Range("J89").Value = 9313.666999999999
$ws.Range("K89").Value = 10043.333
$ws.Range("L89").Value = 46568.335
$ws.Range("M89").Value = -4427.333000000001
$ws.Range("N89").Value = -57800.335
$ws.Range("H91").Value = 5494.0835
$ws.Range("J91").Value = 5286.5713
$ws.Range("L91").Value = 5286.5713
$ws.Range("N91").Value = -8094.5713
$ws.Range("H100").Value = 30636.777
$ws.Range("I100").Value = 34984.805
$ws.Range("J100").Value = 3679
$ws.Range("K100").Value = 34984.805
$ws.Range("L100").Value = 3679
$ws.Range("M100").Value = -34443.805
$ws.Range("N100").Value = -4761
$ws.Range("H111").Value = 2098.3333
$ws.Range("I111").Value = 2098.3333
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 6294.999899999999
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -3227.999899999999
$ws.Range("N111").ClearContents()
$ws.Range("H132").Value = 3063.75
$ws.Range("I132").Value = 3150.8235
$ws.Range("K132").Value = 9452.470499999999
$ws.Range("M132").Value = -6922.470499999999
$ws.Range("H137").Value = 14512.952
$ws.Range("I137").Value = 7668.6
$ws.Range("K137").Value = 23005.8
$ws.Range("M137").Value = -20455.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 188.625
$ws.Range("I5").Value = 151.28572
$ws.Range("K5").Value = 151.28572
$ws.Range("M5").Value = -39.28572
$ws.Range("H12").Value = 3577.5
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H32").Value = 3148.957
$ws.Range("I32").Value = 2292.5967
$ws.Range("K32").Value = 2292.5967
$ws.Range("M32").Value = -2005.5967
$ws.Range("H61").Value = 5429.973
$ws.Range("I61").Value = 4631.913
$ws.Range("J61").Value = 6741.0713
$ws.Range("K61").Value = 4631.913
$ws.Range("L61").Value = 6741.0713
$ws.Range("M61").Value = -4419.913
$ws.Range("N61").Value = -7165.0713
$ws.Range("H74").Value = 3532.4138
$ws.Range("I74").Value = 2907.8235
$ws.Range("J74").Value = 4417.25
$ws.Range("K74").Value = 2907.8235
$ws.Range("L74").Value = 4417.25
$ws.Range("M74").Value = -2033.8235
$ws.Range("N74").Value = -6165.25
$ws.Range("H77").Value = 3532.4138
$ws.Range("I77").Value = 2907.8235
$ws.Range("J77").Value = 4417.25
$ws.Range("K77").Value = 14539.1175
$ws.Range("L77").Value = 22086.25
$ws.Range("M77").Value = -10171.1175
$ws.Range("N77").Value = -30822.25
$ws.Range("H102").Value = 4111.1113
$ws.Range("J102").Value = 3758.8
$ws.Range("L102").Value = 3758.8
$ws.Range("N102").Value = -7002.8
$ws.Range("H110").Value = 771
$ws.Range("I110").Value = 703.5454999999999
$ws.Range("J110").Value = 956.5
$ws.Range("K110").Value = 703.5454999999999
$ws.Range("L110").Value = 956.5
$ws.Range("M110").Value = 1341.4545
$ws.Range("N110").Value = -5046.5
$ws.Range("H122").Value = 1707.4
$ws.Range("I122").Value = 1509.5
$ws.Range("K122").Value = 4528.5
$ws.Range("M122").Value = -2078.5
$ws.Range("H132").Value = 4132.647
$ws.Range("I132").Value = 4132.647
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12397.941
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9867.940999999999
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 5429.973
$ws.Range("I136").Value = 4631.913
$ws.Range("J136").Value = 6741.0713
$ws.Range("K136").Value = 13895.739
$ws.Range("L136").Value = 20223.2139
$ws.Range("M136").Value = -11345.739
$ws.Range("N136").Value = -25323.2139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 188.625
$ws.Range("I4").Value = 151.28572
$ws.Range("K4").Value = 151.28572
$ws.Range("M4").Value = -36.28572
$ws.Range("H20").Value = 26320.4
$ws.Range("I20").Value = 35777.555
$ws.Range("J20").Value = 12134.667
$ws.Range("K20").Value = 35777.555
$ws.Range("L20").Value = 12134.667
$ws.Range("M20").Value = -35530.555
$ws.Range("N20").Value = -12628.667
$ws.Range("H105").Value = 3289.3684
$ws.Range("I105").Value = 3708
$ws.Range("J105").Value = 2382.3333
$ws.Range("K105").Value = 3708
$ws.Range("L105").Value = 2382.3333
$ws.Range("M105").Value = -1961
$ws.Range("N105").Value = -5876.3333
$ws.Range("H107").Value = 1627.6666
$ws.Range("I107").Value = 1609.1
$ws.Range("K107").Value = 1609.1
$ws.Range("M107").Value = 310.9000000000001
$ws.Range("H134").Value = 6166.6904
$ws.Range("I134").Value = 4157.625
$ws.Range("J134").Value = 12595.7
$ws.Range("K134").Value = 12472.875
$ws.Range("L134").Value = 37787.10000000001
$ws.Range("M134").Value = -9937.875
$ws.Range("N134").Value = -42857.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6966.6665
$ws.Range("I6").Value = 19000
$ws.Range("J6").Value = 950
$ws.Range("K6").Value = 19000
$ws.Range("L6").Value = 950
$ws.Range("M6").Value = -18887
$ws.Range("N6").Value = -1176
$ws.Range("H7").Value = 153.3077
$ws.Range("J7").Value = 158.4
$ws.Range("L7").Value = 158.4
$ws.Range("N7").Value = -384.4
$ws.Range("H16").Value = 6295.2144
$ws.Range("I16").Value = 6257.636
$ws.Range("K16").Value = 6257.636
$ws.Range("M16").Value = -5970.636
$ws.Range("H22").Value = 282.10526
$ws.Range("J22").Value = 323.33334
$ws.Range("L22").Value = 323.33334
$ws.Range("N22").Value = -1023.33334
$ws.Range("H31").Value = 2374.5
$ws.Range("I31").Value = 1299.6666
$ws.Range("K31").Value = 1299.6666
$ws.Range("M31").Value = -1004.6666
$ws.Range("H34").Value = 2374.5
$ws.Range("I34").Value = 1299.6666
$ws.Range("K34").Value = 1299.6666
$ws.Range("M34").Value = -1097.6666
$ws.Range("H58").Value = 3584.3044
$ws.Range("I58").Value = 1926.0588
$ws.Range("K58").Value = 1926.0588
$ws.Range("M58").Value = -1723.0588
$ws.Range("H86").Value = 4999
$ws.Range("J86").Value = 4997.5
$ws.Range("L86").Value = 4997.5
$ws.Range("N86").Value = -7243.5
$ws.Range("H89").Value = 4999
$ws.Range("J89").Value = 4997.5
$ws.Range("L89").Value = 24987.5
$ws.Range("N89").Value = -36219.5
$ws.Range("H113").Value = 6295.2144
$ws.Range("I113").Value = 6257.636
$ws.Range("K113").Value = 6257.636
$ws.Range("M113").Value = -4087.636
$ws.Range("H136").Value = 3584.3044
$ws.Range("I136").Value = 1926.0588
$ws.Range("K136").Value = 5778.1764
$ws.Range("M136").Value = -3228.1764
$ws.Range("H141").Value = 129600.336
$ws.Range("I141").Value = 60000
$ws.Range("J141").Value = 135927.64
$ws.Range("K141").Value = 60000
$ws.Range("L141").Value = 135927.64
$ws.Range("M141").Value = -54820
$ws.Range("N141").Value = -146287.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7749.75
$ws.Range("I3").Value = 7749.75
$ws.Range("K3").Value = 23249.25
$ws.Range("M3").Value = -23137.25
$ws.Range("H20").Value = 97.5
$ws.Range("J20").Value = 95
$ws.Range("L20").Value = 285
$ws.Range("N20").Value = -739
$ws.Range("H23").Value = 231.05
$ws.Range("J23").Value = 313.15384
$ws.Range("L23").Value = 939.4615200000001
$ws.Range("N23").Value = -1409.46152
$ws.Range("H98").Value = 711.2222
$ws.Range("J98").Value = 683
$ws.Range("L98").Value = 2049
$ws.Range("N98").Value = -5045
$ws.Range("H107").Value = 855.65
$ws.Range("J107").Value = 989.55554
$ws.Range("L107").Value = 2968.66662
$ws.Range("N107").Value = -6808.66662
$ws.Range("H113").Value = 1267.8334
$ws.Range("I113").Value = 1548.6666
$ws.Range("J113").Value = 987
$ws.Range("K113").Value = 4645.9998
$ws.Range("L113").Value = 2961
$ws.Range("M113").Value = -2475.9998
$ws.Range("N113").Value = -7301
$ws.Range("H132").Value = 10002569
$ws.Range("J132").Value = 20002800
$ws.Range("L132").Value = 180025200
$ws.Range("N132").Value = -180030260
$ws.Range("H133").Value = 1236.75
$ws.Range("I133").Value = 1236.75
$ws.Range("K133").Value = 3710.25
$ws.Range("M133").Value = 1349.75
$ws.Range("H134").Value = 1181.1666
$ws.Range("I134").Value = 1181.1666
$ws.Range("K134").Value = 3543.4998
$ws.Range("M134").Value = 1526.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 81.875
$ws.Range("I2").Value = 39.2
$ws.Range("J2").Value = 153
$ws.Range("K2").Value = 39.2
$ws.Range("L2").Value = 153
$ws.Range("M2").Value = 73.8
$ws.Range("N2").Value = -379
$ws.Range("H5").Value = 18795.6
$ws.Range("I5").Value = 18999.5
$ws.Range("K5").Value = 18999.5
$ws.Range("M5").Value = -18887.5
$ws.Range("H113").Value = 1970.0952
$ws.Range("I113").Value = 1658
$ws.Range("K113").Value = 1658
$ws.Range("M113").Value = 512
$ws.Range("H122").Value = 1608.75
$ws.Range("I122").Value = 1624.2858
$ws.Range("K122").Value = 4872.857400000001
$ws.Range("M122").Value = -2422.857400000001
$ws.Range("H126").Value = 4981
$ws.Range("I126").Value = 4777.2
$ws.Range("K126").Value = 14331.6
$ws.Range("M126").Value = -11861.6
$ws.Range("H132").Value = 3011.5625
$ws.Range("I132").Value = 3084.1482
$ws.Range("J132").Value = 2619.6
$ws.Range("K132").Value = 9252.444600000001
$ws.Range("L132").Value = 7858.799999999999
$ws.Range("M132").Value = -6722.444600000001
$ws.Range("N132").Value = -12918.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 149999
$ws.Range("J17").Value = 149998
$ws.Range("L17").Value = 149998
$ws.Range("N17").Value = -150338
$ws.Range("H22").Value = 1751.4117
$ws.Range("I22").Value = 1647.625
$ws.Range("J22").Value = 1843.6666
$ws.Range("K22").Value = 1647.625
$ws.Range("L22").Value = 1843.6666
$ws.Range("M22").Value = -1352.625
$ws.Range("N22").Value = -2433.6666
$ws.Range("H27").Value = 1751.4117
$ws.Range("I27").Value = 1647.625
$ws.Range("J27").Value = 1843.6666
$ws.Range("K27").Value = 1647.625
$ws.Range("L27").Value = 1843.6666
$ws.Range("M27").Value = -1540.625
$ws.Range("N27").Value = -2057.6666
$ws.Range("H46").Value = 1112.2632
$ws.Range("I46").Value = 935.7692
$ws.Range("J46").Value = 1204.04
$ws.Range("K46").Value = 935.7692
$ws.Range("L46").Value = 1204.04
$ws.Range("M46").Value = -747.7692
$ws.Range("N46").Value = -1580.04
$ws.Range("H55").Value = 90.25
$ws.Range("I55").Value = 84.5
$ws.Range("K55").Value = 84.5
$ws.Range("M55").Value = 88.5
$ws.Range("H61").Value = 2078.261
$ws.Range("I61").Value = 1880
$ws.Range("K61").Value = 1880
$ws.Range("M61").Value = -1678
$ws.Range("H74").Value = 35804.715
$ws.Range("I74").Value = 22584
$ws.Range("J74").Value = 38008.168
$ws.Range("K74").Value = 22584
$ws.Range("L74").Value = 38008.168
$ws.Range("M74").Value = -21586
$ws.Range("N74").Value = -40004.168
$ws.Range("H77").Value = 35804.715
$ws.Range("I77").Value = 22584
$ws.Range("J77").Value = 38008.168
$ws.Range("K77").Value = 67752
$ws.Range("L77").Value = 114024.504
$ws.Range("M77").Value = -62760
$ws.Range("N77").Value = -124008.504
$ws.Range("H113").Value = 2078.261
$ws.Range("I113").Value = 1880
$ws.Range("K113").Value = 1880
$ws.Range("M113").Value = 290
$ws.Range("H132").Value = 2741.92
$ws.Range("I132").Value = 2478.2173
$ws.Range("K132").Value = 7434.651899999999
$ws.Range("M132").Value = -4904.651899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 9000
$ws.Range("I9").Value = 9000
$ws.Range("K9").Value = 9000
$ws.Range("M9").Value = -8860
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H40").Value = 14750
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 14750
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 14750
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -15048
$ws.Range("H41").Value = 30000
$ws.Range("I41").Value = 30000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 30000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -29610
$ws.Range("N41").ClearContents()
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H81").Value = 2254.6667
$ws.Range("I81").Value = 2254.6667
$ws.Range("K81").Value = 4509.3334
$ws.Range("M81").Value = -3448.3334
$ws.Range("H84").Value = 2254.6667
$ws.Range("I84").Value = 2254.6667
$ws.Range("K84").Value = 22546.667
$ws.Range("M84").Value = -17242.667
$ws.Range("H107").Value = 2060.6538
$ws.Range("I107").Value = 1943.25
$ws.Range("J107").Value = 2248.5
$ws.Range("K107").Value = 5829.75
$ws.Range("L107").Value = 6745.5
$ws.Range("M107").Value = -3909.75
$ws.Range("N107").Value = -10585.5
$ws.Range("H113").Value = 138.25
$ws.Range("I113").Value = 270
$ws.Range("J113").Value = 94.333336
$ws.Range("K113").Value = 810
$ws.Range("L113").Value = 283.000008
$ws.Range("M113").Value = 1360
$ws.Range("N113").Value = -4623.000008
$ws.Range("H126").Value = 9800.286
$ws.Range("I126").Value = 13150.5
$ws.Range("J126").Value = 5333.3335
$ws.Range("K126").Value = 39451.5
$ws.Range("L126").Value = 16000.0005
$ws.Range("M126").Value = -36981.5
$ws.Range("N126").Value = -20940.0005
